$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4680
$ws1.Range("F6").Value = 1854
$ws1.Range("F11").Value = 418
$ws1.Range("F13").Value = 1591
$ws1.Range("F14").Value = 819
$ws1.Range("F15").Value = 1819
$ws1.Range("F16").Value = 560
$ws1.Range("F19").Value = 193
$ws1.Range("F20").Value = 1555
$ws1.Range("F26").Value = 1585
$ws1.Range("F31").Value = 4287

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 4680
$ws4.Range("F12").Value = 1854
$ws4.Range("F19").Value = 418
$ws4.Range("F21").Value = 1591
$ws4.Range("F23").Value = 819
$ws4.Range("F24").Value = 1820
$ws4.Range("F25").Value = 560
$ws4.Range("F28").Value = 193
$ws4.Range("F32").Value = 1555
$ws4.Range("F43").Value = 1585
$ws4.Range("F48").Value = 4287
